$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col4a2"
$ws.Cells.Item(2, 3).Value = "Cd93"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 148.2038823333333
$ws.Cells.Item(2, 8).Value = 444.6116469999999
$ws.Cells.Item(2, 9).Value = 0.4690635672357343
$ws.Cells.Item(2, 10).Value = 0.4690635672357342
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 135.955556
$ws.Cells.Item(2, 14).Value = 407.866668
$ws.Cells.Item(2, 15).Value = 0.6947679994035034
$ws.Cells.Item(2, 16).Value = 0.6947679994035034
$ws.Cells.Item(2, 17).Value = 20149.14122398691
$ws.Cells.Item(2, 18).Value = 181342.2710158822
$ws.Cells.Item(2, 19).Value = 0.3258903562014418
$ws.Cells.Item(2, 20).Value = 0.3258903562014418

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col4a2"
$ws.Cells.Item(3, 3).Value = "Cd93"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 148.2038823333333
$ws.Cells.Item(3, 8).Value = 444.6116469999999
$ws.Cells.Item(3, 9).Value = 0.4690635672357343
$ws.Cells.Item(3, 10).Value = 0.4690635672357342
$ws.Cells.Item(3, 11).Value = 2.0
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.449122
$ws.Cells.Item(3, 14).Value = 1.347366
$ws.Cells.Item(3, 15).Value = 0.002295129398228494
$ws.Cells.Item(3, 16).Value = 0.002295129398228494
$ws.Cells.Item(3, 17).Value = 66.56162404131133
$ws.Cells.Item(3, 18).Value = 599.054616371802
$ws.Cells.Item(3, 19).Value = 0.001076561582800661
$ws.Cells.Item(3, 20).Value = 0.001076561582800661

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col4a2"
$ws.Cells.Item(4, 3).Value = "Cd93"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 148.2038823333333
$ws.Cells.Item(4, 8).Value = 444.6116469999999
$ws.Cells.Item(4, 9).Value = 0.4690635672357343
$ws.Cells.Item(4, 10).Value = 0.4690635672357342
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 56.38366533333333
$ws.Cells.Item(4, 14).Value = 169.150996
$ws.Cells.Item(4, 15).Value = 0.2881350899898248
$ws.Cells.Item(4, 16).Value = 0.2881350899898248
$ws.Cells.Item(4, 17).Value = 8356.27810258338
$ws.Cells.Item(4, 18).Value = 75206.5029232504
$ws.Cells.Item(4, 19).Value = 0.1351536731564165
$ws.Cells.Item(4, 20).Value = 0.1351536731564165

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Col4a2"
$ws.Cells.Item(5, 3).Value = "Cd93"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 148.2038823333333
$ws.Cells.Item(5, 8).Value = 444.6116469999999
$ws.Cells.Item(5, 9).Value = 0.4690635672357343
$ws.Cells.Item(5, 10).Value = 0.4690635672357342
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 2.896484
$ws.Cells.Item(5, 14).Value = 8.689452
$ws.Cells.Item(5, 15).Value = 0.01480178120844327
$ws.Cells.Item(5, 16).Value = 0.01480178120844327
$ws.Cells.Item(5, 17).Value = 429.2701739163826
$ws.Cells.Item(5, 18).Value = 3863.431565247443
$ws.Cells.Item(5, 19).Value = 0.006942976295075259
$ws.Cells.Item(5, 20).Value = 0.006942976295075259

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col4a2"
$ws.Cells.Item(6, 3).Value = "Cd93"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 92.51705933333334
$ws.Cells.Item(6, 8).Value = 277.551178
$ws.Cells.Item(6, 9).Value = 0.2928154188528495
$ws.Cells.Item(6, 10).Value = 0.2928154188528495
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 135.955556
$ws.Cells.Item(6, 14).Value = 407.866668
$ws.Cells.Item(6, 15).Value = 0.6947679994035034
$ws.Cells.Item(6, 16).Value = 0.6947679994035034
$ws.Cells.Item(6, 17).Value = 12578.20824114832
$ws.Cells.Item(6, 18).Value = 113203.8741703349
$ws.Cells.Item(6, 19).Value = 0.2034387827508931
$ws.Cells.Item(6, 20).Value = 0.2034387827508931

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col4a2"
$ws.Cells.Item(7, 3).Value = "Cd93"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 92.51705933333334
$ws.Cells.Item(7, 8).Value = 277.551178
$ws.Cells.Item(7, 9).Value = 0.2928154188528495
$ws.Cells.Item(7, 10).Value = 0.2928154188528495
$ws.Cells.Item(7, 11).Value = 2.0
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.449122
$ws.Cells.Item(7, 14).Value = 1.347366
$ws.Cells.Item(7, 15).Value = 0.002295129398228494
$ws.Cells.Item(7, 16).Value = 0.002295129398228494
$ws.Cells.Item(7, 17).Value = 41.55144672190534
$ws.Cells.Item(7, 18).Value = 373.963020497148
$ws.Cells.Item(7, 19).Value = 0.0006720492760637648
$ws.Cells.Item(7, 20).Value = 0.0006720492760637646

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Col4a2"
$ws.Cells.Item(8, 3).Value = "Cd93"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 92.51705933333334
$ws.Cells.Item(8, 8).Value = 277.551178
$ws.Cells.Item(8, 9).Value = 0.2928154188528495
$ws.Cells.Item(8, 10).Value = 0.2928154188528495
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 56.38366533333333
$ws.Cells.Item(8, 14).Value = 169.150996
$ws.Cells.Item(8, 15).Value = 0.2881350899898248
$ws.Cells.Item(8, 16).Value = 0.2881350899898248
$ws.Cells.Item(8, 17).Value = 5216.450911074809
$ws.Cells.Item(8, 18).Value = 46948.05819967329
$ws.Cells.Item(8, 19).Value = 0.08437039706157404
$ws.Cells.Item(8, 20).Value = 0.08437039706157401

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Col4a2"
$ws.Cells.Item(9, 3).Value = "Cd93"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 92.51705933333334
$ws.Cells.Item(9, 8).Value = 277.551178
$ws.Cells.Item(9, 9).Value = 0.2928154188528495
$ws.Cells.Item(9, 10).Value = 0.2928154188528495
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 2.896484
$ws.Cells.Item(9, 14).Value = 8.689452
$ws.Cells.Item(9, 15).Value = 0.01480178120844327
$ws.Cells.Item(9, 16).Value = 0.01480178120844327
$ws.Cells.Item(9, 17).Value = 267.9741820860506
$ws.Cells.Item(9, 18).Value = 2411.767638774456
$ws.Cells.Item(9, 19).Value = 0.004334189764318553
$ws.Cells.Item(9, 20).Value = 0.004334189764318553

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Col4a2"
$ws.Cells.Item(10, 3).Value = "Cd93"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 0.9606106666666666
$ws.Cells.Item(10, 8).Value = 2.881832
$ws.Cells.Item(10, 9).Value = 0.003040321609240458
$ws.Cells.Item(10, 10).Value = 0.003040321609240458
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 135.955556
$ws.Cells.Item(10, 14).Value = 407.866668
$ws.Cells.Item(10, 15).Value = 0.6947679994035034
$ws.Cells.Item(10, 16).Value = 0.6947679994035034
$ws.Cells.Item(10, 17).Value = 130.6003572861973
$ws.Cells.Item(10, 18).Value = 1175.403215575776
$ws.Cells.Item(10, 19).Value = 0.002112318161995233
$ws.Cells.Item(10, 20).Value = 0.002112318161995233

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Col4a2"
$ws.Cells.Item(11, 3).Value = "Cd93"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 0.9606106666666666
$ws.Cells.Item(11, 8).Value = 2.881832
$ws.Cells.Item(11, 9).Value = 0.003040321609240458
$ws.Cells.Item(11, 10).Value = 0.003040321609240458
$ws.Cells.Item(11, 11).Value = 2.0
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.449122
$ws.Cells.Item(11, 14).Value = 1.347366
$ws.Cells.Item(11, 15).Value = 0.002295129398228494
$ws.Cells.Item(11, 16).Value = 0.002295129398228494
$ws.Cells.Item(11, 17).Value = 0.4314313838346667
$ws.Cells.Item(11, 18).Value = 3.882882454512
$ws.Cells.Item(11, 19).Value = 0.000006977931505437139
$ws.Cells.Item(11, 20).Value = 0.000006977931505437137

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Col4a2"
$ws.Cells.Item(12, 3).Value = "Cd93"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 0.9606106666666666
$ws.Cells.Item(12, 8).Value = 2.881832
$ws.Cells.Item(12, 9).Value = 0.003040321609240458
$ws.Cells.Item(12, 10).Value = 0.003040321609240458
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 56.38366533333333
$ws.Cells.Item(12, 14).Value = 169.150996
$ws.Cells.Item(12, 15).Value = 0.2881350899898248
$ws.Cells.Item(12, 16).Value = 0.2881350899898248
$ws.Cells.Item(12, 17).Value = 54.16275034496356
$ws.Cells.Item(12, 18).Value = 487.464753104672
$ws.Cells.Item(12, 19).Value = 0.0008760233404765085
$ws.Cells.Item(12, 20).Value = 0.0008760233404765083

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Col4a2"
$ws.Cells.Item(13, 3).Value = "Cd93"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 0.9606106666666666
$ws.Cells.Item(13, 8).Value = 2.881832
$ws.Cells.Item(13, 9).Value = 0.003040321609240458
$ws.Cells.Item(13, 10).Value = 0.003040321609240458
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 2.896484
$ws.Cells.Item(13, 14).Value = 8.689452
$ws.Cells.Item(13, 15).Value = 0.01480178120844327
$ws.Cells.Item(13, 16).Value = 0.01480178120844327
$ws.Cells.Item(13, 17).Value = 2.782393426229333
$ws.Cells.Item(13, 18).Value = 25.04154083606399
$ws.Cells.Item(13, 19).Value = 0.00004500217526327943
$ws.Cells.Item(13, 20).Value = 0.00004500217526327942

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Col4a2"
$ws.Cells.Item(14, 3).Value = "Cd93"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3.0
$ws.Cells.Item(14, 6).Value = 1.0
$ws.Cells.Item(14, 7).Value = 74.27537266666666
$ws.Cells.Item(14, 8).Value = 222.826118
$ws.Cells.Item(14, 9).Value = 0.2350806923021759
$ws.Cells.Item(14, 10).Value = 0.2350806923021759
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 12).Value = 1.0
$ws.Cells.Item(14, 13).Value = 135.955556
$ws.Cells.Item(14, 14).Value = 407.866668
$ws.Cells.Item(14, 15).Value = 0.6947679994035034
$ws.Cells.Item(14, 16).Value = 0.6947679994035034
$ws.Cells.Item(14, 17).Value = 10098.14958800387
$ws.Cells.Item(14, 18).Value = 90883.34629203481
$ws.Cells.Item(14, 19).Value = 0.1633265422891733
$ws.Cells.Item(14, 20).Value = 0.1633265422891733

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Col4a2"
$ws.Cells.Item(15, 3).Value = "Cd93"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3.0
$ws.Cells.Item(15, 6).Value = 1.0
$ws.Cells.Item(15, 7).Value = 74.27537266666666
$ws.Cells.Item(15, 8).Value = 222.826118
$ws.Cells.Item(15, 9).Value = 0.2350806923021759
$ws.Cells.Item(15, 10).Value = 0.2350806923021759
$ws.Cells.Item(15, 11).Value = 2.0
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.449122
$ws.Cells.Item(15, 14).Value = 1.347366
$ws.Cells.Item(15, 15).Value = 0.002295129398228494
$ws.Cells.Item(15, 16).Value = 0.002295129398228494
$ws.Cells.Item(15, 17).Value = 33.35870392279866
$ws.Cells.Item(15, 18).Value = 300.228335305188
$ws.Cells.Item(15, 19).Value = 0.0005395406078586308
$ws.Cells.Item(15, 20).Value = 0.0005395406078586306

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Col4a2"
$ws.Cells.Item(16, 3).Value = "Cd93"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3.0
$ws.Cells.Item(16, 6).Value = 1.0
$ws.Cells.Item(16, 7).Value = 74.27537266666666
$ws.Cells.Item(16, 8).Value = 222.826118
$ws.Cells.Item(16, 9).Value = 0.2350806923021759
$ws.Cells.Item(16, 10).Value = 0.2350806923021759
$ws.Cells.Item(16, 11).Value = 3.0
$ws.Cells.Item(16, 12).Value = 1.0
$ws.Cells.Item(16, 13).Value = 56.38366533333333
$ws.Cells.Item(16, 14).Value = 169.150996
$ws.Cells.Item(16, 15).Value = 0.2881350899898248
$ws.Cells.Item(16, 16).Value = 0.2881350899898248
$ws.Cells.Item(16, 17).Value = 4187.917754945947
$ws.Cells.Item(16, 18).Value = 37691.25979451352
$ws.Cells.Item(16, 19).Value = 0.06773499643135777
$ws.Cells.Item(16, 20).Value = 0.06773499643135776

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Col4a2"
$ws.Cells.Item(17, 3).Value = "Cd93"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3.0
$ws.Cells.Item(17, 6).Value = 1.0
$ws.Cells.Item(17, 7).Value = 74.27537266666666
$ws.Cells.Item(17, 8).Value = 222.826118
$ws.Cells.Item(17, 9).Value = 0.2350806923021759
$ws.Cells.Item(17, 10).Value = 0.2350806923021759
$ws.Cells.Item(17, 11).Value = 3.0
$ws.Cells.Item(17, 12).Value = 1.0
$ws.Cells.Item(17, 13).Value = 2.896484
$ws.Cells.Item(17, 14).Value = 8.689452
$ws.Cells.Item(17, 15).Value = 0.01480178120844327
$ws.Cells.Item(17, 16).Value = 0.01480178120844327
$ws.Cells.Item(17, 17).Value = 215.1374285230373
$ws.Cells.Item(17, 18).Value = 1936.236856707336
$ws.Cells.Item(17, 19).Value = 0.003479612973786182
$ws.Cells.Item(17, 20).Value = 0.003479612973786182
